$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row right after the header so the "Not applicable" (-1) entry
# becomes the first lookup value, shifting every other row down by one.
$ws.Rows.Item(2).Insert()

$ws.Range("A2").Value = -1
$ws.Range("B2").Value = "Not applicable"

# Keep the named range in sync with the table's new extent (one extra row).
$wb.Names.Item("dbo_prepmet").RefersTo = "=dbo_prepmet!`$A`$1:`$G`$16"

# Match the author's final cursor position.
$ws.Range("A3").Select()
